# Commit: "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# The "mora" (overdue) detail table on Hoja1 (rows 16-36, columns C:G) is
# rebuilt: the old periods (1906-1912) for each of the three workers are
# replaced with a new arrangement where, for every worker, all seven
# periods (1912 down to 1906) are listed together in consecutive rows,
# along with the corresponding "Valor Mora" (F) and "Salario Basico" (G)
# amounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: DocTrabajador, NombreTrabajador, PeriodoMora, ValorMora, SalarioBasico
$data = @(
    @("73147208", "GREGORIO JOSE DIAZ BANDERA", "1912", 34133, 1600000),
    @("73147208", "GREGORIO JOSE DIAZ BANDERA", "1911", 64000, 1600000),
    @("73147208", "GREGORIO JOSE DIAZ BANDERA", "1910", 64000, 1600000),
    @("73147208", "GREGORIO JOSE DIAZ BANDERA", "1909", 64000, 1600000),
    @("73147208", "GREGORIO JOSE DIAZ BANDERA", "1908", 64000, 1600000),
    @("73147208", "GREGORIO JOSE DIAZ BANDERA", "1907", 64000, 1600000),
    @("73147208", "GREGORIO JOSE DIAZ BANDERA", "1906", 64000, 1600000),
    @("1047453567", "RUBEN REYES MUÑOZ", "1912", 17667, 828116),
    @("1047453567", "RUBEN REYES MUÑOZ", "1911", 33125, 828116),
    @("1047453567", "RUBEN REYES MUÑOZ", "1910", 33125, 828116),
    @("1047453567", "RUBEN REYES MUÑOZ", "1909", 33125, 828116),
    @("1047453567", "RUBEN REYES MUÑOZ", "1908", 33125, 828116),
    @("1047453567", "RUBEN REYES MUÑOZ", "1907", 33125, 828116),
    @("1047453567", "RUBEN REYES MUÑOZ", "1906", 33125, 828116),
    @("79539228", "OSCAR HUMBERTO CAÑAS DIAZ", "1912", 34133, 1600000),
    @("79539228", "OSCAR HUMBERTO CAÑAS DIAZ", "1911", 64000, 1600000),
    @("79539228", "OSCAR HUMBERTO CAÑAS DIAZ", "1910", 64000, 1600000),
    @("79539228", "OSCAR HUMBERTO CAÑAS DIAZ", "1909", 64000, 1600000),
    @("79539228", "OSCAR HUMBERTO CAÑAS DIAZ", "1908", 64000, 1600000),
    @("79539228", "OSCAR HUMBERTO CAÑAS DIAZ", "1907", 64000, 1600000),
    @("79539228", "OSCAR HUMBERTO CAÑAS DIAZ", "1906", 64000, 1600000)
)

$startRow = 16
$r = $startRow
foreach ($row in $data) {
    $ws.Cells.Item($r, 3).Value = $row[0]   # C: N° Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $row[1]   # D: Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $row[2]   # E: Periodo Mora
    $ws.Cells.Item($r, 6).Value = $row[3]   # F: Valor Mora
    $ws.Cells.Item($r, 7).Value = $row[4]   # G: Salario Basico
    $r++
}
